# Generate Report for Handoff
# Updates the localization-status report to reflect a freshly generated
# handoff package: the "Status" text flips from the old handback message to
# "Ready for handoff", and the associated timestamps move forward a little
# over a minute. The Status column also narrows to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-16 15:00:47"

# Status columns (E & F) shrink to fit the shorter "Ready for handoff" text.
$overview.Range("E1").ColumnWidth = 16.29
$overview.Range("F1").ColumnWidth = 16.29

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-16 15:00:41"
$zhcn.Range("C1").ColumnWidth = 16.29

# --- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-16 15:00:47"
$dede.Range("C1").ColumnWidth = 16.29
